$d = $word.ActiveDocument

# 1. DATE field cached text: February 11, 2010 -> February 12, 2010
$d.Content.Find.Execute("February 11, 2010", $true, $false, $false, $false, $false, $true, 1, $false, "February 12, 2010", 2) | Out-Null

# 2. Getting-started list rework.
# Remove the hyperlink (converts it to plain run text) so the wildcard
# replace below can swallow the whole paragraph (incl. the old URL) in one shot.
$d.Hyperlinks.Item(1).Delete()

$p19 = $d.Paragraphs.Item(19)
$p19.Range.Find.Execute("Connect to*increment over time.", $false, $false, $true, $false, $false, $true, 1, $false, "If you haven’t already, Clone the Google Code repository to some repository location on your local disk. This will get the baseline database and all other PRGFX code.", 2) | Out-Null

# "Alternately, if you are on the City network, you may:" ilvl 1 -> 0 (text unchanged)
$p20 = $d.Paragraphs.Item(20)
$p20.Range.ListFormat.ListLevelNumber = 1

# "Create an empty EA project..." ilvl 2 -> 1
$p21 = $d.Paragraphs.Item(21)
$p21.Range.ListFormat.ListLevelNumber = 2

# "Run Tools > Data Management > Project Transfer..." ilvl 2 -> 1
$p22 = $d.Paragraphs.Item(22)
$p22.Range.ListFormat.ListLevelNumber = 2

# Insert a new non-list paragraph after it (before "Make a copy...")
$p23target = $d.Paragraphs.Item(23)
$p23target.Range.InsertParagraphBefore()
$newp = $d.Paragraphs.Item(23)
$newp.Range.Text = "This won’t get you the rest of the code you need, so you will want to do step 1 in any case."
$newp.Range.ListFormat.RemoveNumbers()
$newp.Range.ParagraphFormat.LeftIndent = 36

# "Make a copy of this file..." gets a new trailing sentence
$p24 = $d.Paragraphs.Item(24)
$p24.Range.Find.Execute("Make a copy of this file to some other working area.", $false, $false, $false, $false, $false, $true, 1, $false, "Make a copy of this file to some other working area. You don’t want to work directly in your controlled area, as you might inadvertently overwrite the baseline EA database.", 2) | Out-Null

# "Using TortoiseHg, clone the EA repository..." -> replaced text
$p26 = $d.Paragraphs.Item(26)
$p26.Range.Find.Execute("Using TortoiseHg, clone the EA repository*working copy.", $false, $false, $true, $false, $false, $true, 1, $false, "The various XML export / import files will be a part of the cloned Hg repository.", 2) | Out-Null

# "At this point you should have an environment that looks something like this figure 1." -> drop "this "
$p27 = $d.Paragraphs.Item(27)
$p27.Range.Find.Execute("like this figure 1.", $false, $false, $false, $false, $false, $true, 1, $false, "like figure 1.", 2) | Out-Null

# 4. "Pull any Google hosted changes" gets a new trailing sentence
$d.Content.Find.Execute("Pull any Google hosted changes", $true, $false, $false, $false, $false, $true, 1, $false, "Pull any Google hosted changes. The ‘after pull’ action should be “update”.", 2) | Out-Null

# 5. New paragraph near the end, after the trailing blank BodyText paragraph.
$lastP = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastP.Range.InsertParagraphAfter()
$newEndP = $d.Paragraphs.Item($d.Paragraphs.Count)
$newEndP.Range.Text = "These procedures mimic what EA does with its version control mechanism turned on; only we’re doing it by hand."

Write-Host "done"
